# Atualizado por script em 12-11-2023 20:45
#
# This script:
#  1) Swaps the match-detail columns (F:V) between six pairs of rows whose
#     fixture order had been reversed (the "Indice"/date/league columns
#     A:E stay put, only the home/away/odds/url details swap).
#  2) Appends four new fixture rows (176-179) at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Swap columns F..V between the following row pairs
# ---------------------------------------------------------------------------
$swapPairs = @(
    @(62, 63),
    @(64, 65),
    @(80, 81),
    @(84, 85),
    @(102, 103),
    @(136, 137)
)

foreach ($pair in $swapPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    for ($col = 6; $col -le 22; $col++) {
        $c1 = $ws.Cells.Item($r1, $col)
        $c2 = $ws.Cells.Item($r2, $col)
        $v1 = $c1.Value2
        $v2 = $c2.Value2
        $c1.Value = $v2
        $c2.Value = $v1
    }
}

# ---------------------------------------------------------------------------
# 2) Append new rows 176-179
# ---------------------------------------------------------------------------
$newRows = @(
    @{ Row=176; A=175; E=45242.77083333334; F="Colon Santa Fe";     G=3; H="Talleres Cordoba";  I=0; J=2.64; K="07/11/2023 05:42"; L=2.39; M="12/11/2023 18:21"; N=3.15; O="07/11/2023 05:42"; P=3.22; Q="12/11/2023 18:21"; R=2.88; S="07/11/2023 05:42"; T=3.3;  U="12/11/2023 18:21"; V="https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/colon-santa-fe-talleres-cordoba/KpXeyZQg/" },
    @{ Row=177; A=176; E=45242.77083333334; F="Defensa y Justicia"; G=0; H="San Lorenzo";        I=1; J=2.65; K="08/11/2023 23:13"; L=2.73; M="12/11/2023 18:24"; N=2.86; O="08/11/2023 23:13"; P=2.85; Q="12/11/2023 18:24"; R=3.16; S="08/11/2023 23:13"; T=3.17; U="12/11/2023 18:24"; V="https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/defensa-y-justicia-san-lorenzo/C0o5XjQI/" },
    @{ Row=178; A=177; E=45242.86458333334; F="Boca Juniors";       G=1; H="Newells Old Boys";   I=0; J=1.96; K="08/11/2023 23:12"; L=2.06; M="12/11/2023 20:41"; N=3.14; O="08/11/2023 23:12"; P=3.03; Q="12/11/2023 20:41"; R=4.24; S="08/11/2023 23:12"; T=4.57; U="12/11/2023 20:41"; V="https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/boca-juniors-newells-old-boys/Iwn1YWAC/" },
    @{ Row=179; A=178; E=45242.86458333334; F="Tigre";              G=1; H="Platense";           I=1; J=2.03; K="07/11/2023 05:42"; L=1.93; M="12/11/2023 20:43"; N=3.13; O="07/11/2023 05:42"; P=3.21; Q="12/11/2023 20:40"; R=4.31; S="07/11/2023 05:42"; T=4.83; U="12/11/2023 20:43"; V="https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/tigre-platense/byWazgt0/" }
)

$templateRow = 175

foreach ($rowData in $newRows) {
    $r = $rowData.Row

    # Clone the whole template row first: this carries over the exact
    # cell styles (bold/border/centered A column, date-time-formatted E
    # column, unstyled rest) AND the exact cell data types -- in
    # particular column D ("2023") is stored as text in the template row,
    # and Excel would otherwise re-interpret a freshly assigned "2023"
    # string as a number. Column D's value is identical in every new row,
    # so it is intentionally left untouched after the clone.
    $ws.Range("A${templateRow}:V${templateRow}").Copy($ws.Range("A${r}:V${r}"))

    $ws.Range("A$r").Value = $rowData.A
    $ws.Range("B$r").Value = "argentina"
    $ws.Range("C$r").Value = "copa-de-la-liga-profesional"
    # D$r intentionally left as copied ("2023" text)
    $ws.Range("E$r").Value = $rowData.E
    $ws.Range("F$r").Value = $rowData.F
    $ws.Range("G$r").Value = $rowData.G
    $ws.Range("H$r").Value = $rowData.H
    $ws.Range("I$r").Value = $rowData.I
    $ws.Range("J$r").Value = $rowData.J
    $ws.Range("K$r").Value = $rowData.K
    $ws.Range("L$r").Value = $rowData.L
    $ws.Range("M$r").Value = $rowData.M
    $ws.Range("N$r").Value = $rowData.N
    $ws.Range("O$r").Value = $rowData.O
    $ws.Range("P$r").Value = $rowData.P
    $ws.Range("Q$r").Value = $rowData.Q
    $ws.Range("R$r").Value = $rowData.R
    $ws.Range("S$r").Value = $rowData.S
    $ws.Range("T$r").Value = $rowData.T
    $ws.Range("U$r").Value = $rowData.U
    $ws.Range("V$r").Value = $rowData.V
}
